# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de):
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - A new "Latest Target File" (F) hyperlink is recorded, mirroring the
#     source markdown file already linked from column A
#   - A new "Latest Handback File" (G) hyperlink is recorded, mirroring the
#     handed-off xlf file already linked from column D
#   - "Latest Handback DateTime" (H) moves from the zero-date sentinel to the
#     real timestamp the handback report was generated at
#
# The same status-text change also ripples into the "Overview" sheet, which
# shares the very same "Ready for handoff" string in its own Status columns.

function Get-HyperlinkAt($ws, $row, $col) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Row -eq $row -and $h.Range.Column -eq $col) {
            return $h
        }
    }
    return $null
}

function Copy-Hyperlink($ws, $fromRow, $fromCol, $toRow, $toCol) {
    $src = Get-HyperlinkAt $ws $fromRow $fromCol
    $dst = $ws.Cells.Item($toRow, $toCol)
    $dst.Value = $src.TextToDisplay
    $dst.Style = "HyperLink"
    $ws.Hyperlinks.Add($dst, $src.Address, "", "", $src.TextToDisplay) | Out-Null
}

$wb = $excel.ActiveWorkbook
$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: same shared Status text, columns B (zh-cn) and C (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

Copy-Hyperlink $zh 2 1 2 6   # A2 -> F2 (Latest Target File)
Copy-Hyperlink $zh 2 4 2 7   # D2 -> G2 (Latest Handback File)
Copy-Hyperlink $zh 3 1 3 6   # A3 -> F3
Copy-Hyperlink $zh 3 4 3 7   # D3 -> G3

$zh.Range("H2").Value = "2016-03-12 16:38:49"
$zh.Range("H3").Value = "2016-03-12 16:38:49"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

Copy-Hyperlink $de 2 1 2 6   # A2 -> F2
Copy-Hyperlink $de 2 4 2 7   # D2 -> G2
Copy-Hyperlink $de 3 1 3 6   # A3 -> F3
Copy-Hyperlink $de 3 4 3 7   # D3 -> G3

$de.Range("H2").Value = "2016-03-12 16:38:54"
$de.Range("H3").Value = "2016-03-12 16:38:54"
